$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastDataRow = 148
$newLink88 = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-12-august-2025.pdf"

# --- Step 1: capture existing F-column hyperlink target URLs (rows 2..148) before
#             touching anything, keyed by their (pre-insert) row number. The text
#             shown in the cell already equals the hyperlink target. ---
$urls = @{}
for ($r = 2; $r -le $lastDataRow; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $val = $cell.Value()
    if ($val) {
        $urls[$r] = $val
    }
}

# --- Step 2: remove the existing hyperlink objects. Inserting a row shifts cell
#             text/values correctly, but leaves hyperlink anchors pointing at their
#             old row numbers, so we rebuild them afterwards instead. ---
for ($r = 2; $r -le $lastDataRow; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $hl = $cell.Hyperlinks()
    if ($hl.Count() -gt 0) {
        $hl.Delete()
    }
}

# --- Step 3: insert a new row at position 2; old rows 2..148 shift down to 3..149 ---
$ws.Rows.Item(2).Insert()

# --- Step 4: populate the new row 2 with the latest price entry. Columns B-F reuse
#             the values that used to live in row 2 (now duplicated into row 3).
#             Force text format first on the text-ish columns so date-like
#             strings aren't auto-converted to date serials (column D holds a
#             genuine number and must stay numeric). ---
$ws.Cells.Item(2, 1).NumberFormat = "@"
$ws.Cells.Item(2, 2).NumberFormat = "@"
$ws.Cells.Item(2, 3).NumberFormat = "@"
$ws.Cells.Item(2, 5).NumberFormat = "@"
$ws.Cells.Item(2, 6).NumberFormat = "@"

$ws.Cells.Item(2, 1).Value = "06-11-2025"
$ws.Cells.Item(2, 2).Value = $ws.Cells.Item(3, 2).Value()
$ws.Cells.Item(2, 3).Value = $ws.Cells.Item(3, 3).Value()
$ws.Cells.Item(2, 4).Value = $ws.Cells.Item(3, 4).Value()
$ws.Cells.Item(2, 5).Value = $ws.Cells.Item(3, 5).Value()
$ws.Cells.Item(2, 6).Value = $ws.Cells.Item(3, 6).Value()

# --- Step 5: copy cell formatting from row 3 back onto row 2 now that the values
#             are set, so row 2 ends up with the same style (General number
#             format, borders, alignment) as the rest of the table. ---
$ws.Range("A3:F3").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Step 6: re-create all hyperlinks, shifted down one row to match the data,
#             using the URLs captured in step 1. ---
foreach ($r in $urls.Keys) {
    $target = $urls[$r]
    $newRow = $r + 1
    $ws.Hyperlinks().Add($ws.Cells.Item($newRow, 6), $target)
}

# --- Step 7: the diff also adds a brand-new hyperlink on F88 (previously that row
#             had no circular link at all). ---
$ws.Hyperlinks().Add($ws.Cells.Item(88, 6), $newLink88)
